$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A20").Value = "fi_lidar_rain_intensity"
$ws.Range("A19").Value = "fi_lidar_rain_reflectivity"
$ws.Range("A20").Select()
